# country_parameters.xlsx - add Hydro columns, switch price currency label
# to USD, tighten interest-rate number formats, refresh scenario data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two new columns (Hydro interest rate / Hydro lifetime)
#    right after "Wind lifetime (years)" (col G) and before
#    "Plant interest rate" (old col H).
# ---------------------------------------------------------------------
$ws.Columns("H:I").Insert()

# ---------------------------------------------------------------------
# 2. Header row text
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Electricity price (USD/kWh)"
$ws.Range("C1").Value = "Heat price (USD/kWh)"
$ws.Range("D1").Value = "Solar interest rate"
$ws.Range("E1").Value = "Solar lifetime (years)"
$ws.Range("F1").Value = "Wind interest rate"
$ws.Range("G1").Value = "Wind lifetime (years)"
$ws.Range("H1").Value = "Hydro interest rate"
$ws.Range("I1").Value = "Hydro lifetime"
$ws.Range("J1").Value = "Plant interest rate"
$ws.Range("K1").Value = "Plant lifetime (years)"
$ws.Range("L1").Value = "Infrastructure interest rate"
$ws.Range("M1").Value = "Infrastructure lifetime (years)"

# ---------------------------------------------------------------------
# 3. Data rows (Laos / Other) - identical scenario values
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Laos"
$ws.Range("A3").Value = "Other"

foreach ($r in 2,3) {
    $ws.Range("B$r").Value = 0.06
    $ws.Range("C$r").Value = 0.028
    $ws.Range("D$r").Value = 0.1
    $ws.Range("E$r").Value = 20
    $ws.Range("F$r").Value = 0.1
    $ws.Range("G$r").Value = 20
    $ws.Range("H$r").Value = 0.1
    $ws.Range("I$r").Value = 80
    $ws.Range("J$r").Value = 0.1
    $ws.Range("K$r").Value = 20
    $ws.Range("L$r").Value = 0.1
    $ws.Range("M$r").Value = 40
}

# ---------------------------------------------------------------------
# 4. Number formats
#    - Electricity price column keeps 2-decimal currency format
#    - all other rate/price columns switch to 3-decimal format
#    - lifetime/year columns keep the integer format
# ---------------------------------------------------------------------
$ws.Range("B1:B3").NumberFormat = "#,##0.00"
$ws.Range("C1:C3,D1:D3,F1:F3,H1:H3,J1:J3,L1:L3").NumberFormat = "#,##0.000"
$ws.Range("E1:E3,G1:G3,I1:I3,K1:K3,M1:M3").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# 6. Row heights
# ---------------------------------------------------------------------
$ws.Rows("1:3").RowHeight = 19.5

# ---------------------------------------------------------------------
# 7. Column widths (best-fit on new text / new columns)
# ---------------------------------------------------------------------
$ws.Columns("A:M").AutoFit()

Write-Output "done"
